$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table1")
$lo = $ws.ListObjects.Item(1)

# Copy formatting from the last existing data row down across all new rows
$ws.Range("A118:G118").Copy()
$ws.Range("A119:G190").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A119").Value = 45474.43814814815
$ws.Range("B119").Value = 45474.438472222224
$ws.Range("C119").Value = 'Excel tables'
$ws.Range("D119").Value = 45474
$ws.Range("E119").Value = 'Yes'
$ws.Range("F119").Value = 'About right'
$ws.Range("G119").Value = 'Just brilliant as ever. would recommend these to everyone!'

$ws.Range("A120").Value = 45474.43834490741
$ws.Range("B120").Value = 45474.438576388886
$ws.Range("C120").Value = 'Excel tables'
$ws.Range("D120").Value = 45474
$ws.Range("E120").Value = 'Yes'
$ws.Range("F120").Value = 'About right'
$ws.Range("G120").Value = 'Great thans so much really usefull '

$ws.Range("A121").Value = 45474.43818287037
$ws.Range("B121").Value = 45474.43876157407
$ws.Range("C121").Value = 'Excel tables'
$ws.Range("D121").Value = 45474
$ws.Range("E121").Value = 'Yes'
$ws.Range("F121").Value = 'About right'
$ws.Range("G121").Value = 'good introduction'

$ws.Range("A122").Value = 45474.43875
$ws.Range("B122").Value = 45474.43881944445
$ws.Range("C122").Value = 'Excel tables'
$ws.Range("D122").Value = 45474
$ws.Range("E122").Value = 'Yes'

$ws.Range("A123").Value = 45474.43846064815
$ws.Range("B123").Value = 45474.438888888886
$ws.Range("C123").Value = 'Excel tables'
$ws.Range("D123").Value = 45474
$ws.Range("E123").Value = 'Yes'
$ws.Range("F123").Value = 'About right'
$ws.Range("G123").Value = 'Clear, presise guidance & instruction.'

$ws.Range("A124").Value = 45474.43855324074
$ws.Range("B124").Value = 45474.43898148148
$ws.Range("C124").Value = 'Excel tables'
$ws.Range("D124").Value = 45474
$ws.Range("E124").Value = 'Yes'
$ws.Range("F124").Value = 'About right'
$ws.Range("G124").Value = 'Thought it was very engaging '

$ws.Range("A125").Value = 45474.43834490741
$ws.Range("B125").Value = 45474.43902777778
$ws.Range("C125").Value = 'Excel tables'
$ws.Range("D125").Value = 45474
$ws.Range("E125").Value = 'Yes'
$ws.Range("F125").Value = 'About right'
$ws.Range("G125").Value = 'Useful and informative session.  Agree 45mins would be better '

$ws.Range("A126").Value = 45474.43861111111
$ws.Range("B126").Value = 45474.439039351855
$ws.Range("C126").Value = 'Excel tables'
$ws.Range("D126").Value = 45474
$ws.Range("E126").Value = 'Yes'
$ws.Range("F126").Value = 'About right'
$ws.Range("G126").Value = 'Excellent session, very informative'

$ws.Range("A127").Value = 45474.43828703704
$ws.Range("B127").Value = 45474.439155092594
$ws.Range("C127").Value = 'Excel tables'
$ws.Range("D127").Value = 45474
$ws.Range("E127").Value = 'Yes'
$ws.Range("F127").Value = 'About right'
$ws.Range("G127").Value = 'Really enjoyed the session, these are brilliant, thanks Brendan.'

$ws.Range("A128").Value = 45474.43849537037
$ws.Range("B128").Value = 45474.43922453704
$ws.Range("C128").Value = 'Excel tables'
$ws.Range("D128").Value = 45474
$ws.Range("E128").Value = 'Yes'
$ws.Range("F128").Value = 'About right'
$ws.Range("G128").Value = 'Very useful beginner session and easy to follow '

$ws.Range("A129").Value = 45474.43858796296
$ws.Range("B129").Value = 45474.43938657407
$ws.Range("C129").Value = 'Excel tables'
$ws.Range("D129").Value = 45474
$ws.Range("E129").Value = 'Yes'
$ws.Range("F129").Value = 'About right'
$ws.Range("G129").Value = 'Really easy to understand and tables it is from now on :)'

$ws.Range("A130").Value = 45474.438263888886
$ws.Range("B130").Value = 45474.439571759256
$ws.Range("C130").Value = 'Excel tables'
$ws.Range("D130").Value = 45474
$ws.Range("E130").Value = 'Yes'
$ws.Range("F130").Value = 'About right'
$ws.Range("G130").Value = 'Really useful, I have been using excel for many years and I am self taught so really good to get some nuggets that I may have missed along the way and I now understand how important tables are! Off to change all my excel sheets to tables :)'

$ws.Range("A131").Value = 45474.438576388886
$ws.Range("B131").Value = 45474.44001157407
$ws.Range("C131").Value = 'Excel tables'
$ws.Range("D131").Value = 45474
$ws.Range("E131").Value = 'Yes'
$ws.Range("F131").Value = 'About right'

$ws.Range("A132").Value = 45474.43817129629
$ws.Range("B132").Value = 45474.44011574074
$ws.Range("C132").Value = 'Excel tables'
$ws.Range("D132").Value = 45474
$ws.Range("E132").Value = 'Yes'
$ws.Range("F132").Value = 'About right'
$ws.Range("G132").Value = 'Tables came in after I''d learnt Excel & I''d used them occasionally but not fully understood & utilised. Now feel confident to make tables at the start of any Excel work.'

$ws.Range("A133").Value = 45474.43844907408
$ws.Range("B133").Value = 45474.4402662037
$ws.Range("C133").Value = 'Excel tables'
$ws.Range("D133").Value = 45474
$ws.Range("E133").Value = 'Yes'
$ws.Range("F133").Value = 'About right'
$ws.Range("G133").Value = 'Really useful for beginners, nice length (about coffee break in length :) Really enjoyed and off to make some tables!'

$ws.Range("A134").Value = 45474.43965277778
$ws.Range("B134").Value = 45474.44027777778
$ws.Range("C134").Value = 'Excel tables'
$ws.Range("D134").Value = 45474
$ws.Range("E134").Value = 'Yes'
$ws.Range("F134").Value = 'Too easy / simple relative to description'
$ws.Range("G134").Value = 'learned a new way to make tables'

$ws.Range("A135").Value = 45474.438252314816
$ws.Range("B135").Value = 45474.440347222226
$ws.Range("C135").Value = 'Excel tables'
$ws.Range("D135").Value = 45474
$ws.Range("E135").Value = 'Yes'
$ws.Range("F135").Value = 'Too easy / simple relative to description'
$ws.Range("G135").Value = 'Material clear with good structure, and probably good for some folk, however I personally did not learn anything new with this one.  '

$ws.Range("A136").Value = 45474.43863425926
$ws.Range("B136").Value = 45474.44076388889
$ws.Range("C136").Value = 'Excel tables'
$ws.Range("D136").Value = 45474
$ws.Range("E136").Value = 'Yes'
$ws.Range("F136").Value = 'About right'
$ws.Range("G136").Value = 'Understood the actual how to, not quite sure I got why I would use a table, possibly because I don''t have a need for it at the moment. Love the short sessions'

$ws.Range("A137").Value = 45474.43827546296
$ws.Range("B137").Value = 45474.44142361111
$ws.Range("C137").Value = 'Excel tables'
$ws.Range("D137").Value = 45474
$ws.Range("E137").Value = 'Yes'
$ws.Range("F137").Value = 'About right'
$ws.Range("G137").Value = 'Now I know some quick and easy tips for tables and why to use them.'

$ws.Range("A138").Value = 45474.447164351855
$ws.Range("B138").Value = 45474.447488425925
$ws.Range("C138").Value = 'Excel tables'
$ws.Range("D138").Value = 45474
$ws.Range("E138").Value = 'Yes'
$ws.Range("F138").Value = 'About right'
$ws.Range("G138").Value = 'Very well presented and engaging '

$ws.Range("A139").Value = 45474.49364583333
$ws.Range("B139").Value = 45474.4943287037
$ws.Range("C139").Value = 'Tidy data in Excel'
$ws.Range("D139").Value = 45471
$ws.Range("E139").Value = 'Yes'
$ws.Range("F139").Value = 'About right'
$ws.Range("G139").Value = 'Very informative session. '

$ws.Range("A140").Value = 45474.5834375
$ws.Range("B140").Value = 45474.58353009259
$ws.Range("C140").Value = 'An introduction to AI'
$ws.Range("D140").Value = 45474
$ws.Range("E140").Value = 'Yes'
$ws.Range("F140").Value = 'About right'

$ws.Range("A141").Value = 45474.58346064815
$ws.Range("B141").Value = 45474.58361111111
$ws.Range("C141").Value = 'An introduction to AI'
$ws.Range("D141").Value = 45474
$ws.Range("E141").Value = 'Yes'
$ws.Range("F141").Value = 'About right'
$ws.Range("G141").Value = 'Chatty, fun, informative'

$ws.Range("A142").Value = 45474.58347222222
$ws.Range("B142").Value = 45474.583657407406
$ws.Range("C142").Value = 'An introduction to AI'
$ws.Range("D142").Value = 45474
$ws.Range("E142").Value = 'Yes'
$ws.Range("F142").Value = 'About right'

$ws.Range("A143").Value = 45474.583333333336
$ws.Range("B143").Value = 45474.58372685185
$ws.Range("C143").Value = 'An introduction to AI'
$ws.Range("D143").Value = 45474
$ws.Range("E143").Value = 'Yes'
$ws.Range("F143").Value = 'About right'
$ws.Range("G143").Value = 'This was a great session. Thanks'

$ws.Range("A144").Value = 45474.58341435185
$ws.Range("B144").Value = 45474.583761574075
$ws.Range("C144").Value = 'An introduction to AI'
$ws.Range("D144").Value = 45474
$ws.Range("E144").Value = 'Yes'
$ws.Range("F144").Value = 'About right'
$ws.Range("G144").Value = 'thought- stimulating and informative session'

$ws.Range("A145").Value = 45474.583333333336
$ws.Range("B145").Value = 45474.58392361111
$ws.Range("C145").Value = 'An introduction to AI'
$ws.Range("D145").Value = 45474
$ws.Range("E145").Value = 'Yes'
$ws.Range("F145").Value = 'About right'
$ws.Range("G145").Value = 'Important for colleagues to have access to this sort of content. '

$ws.Range("A146").Value = 45474.583402777775
$ws.Range("B146").Value = 45474.58398148148
$ws.Range("C146").Value = 'An introduction to AI'
$ws.Range("D146").Value = 45474
$ws.Range("E146").Value = 'Yes'
$ws.Range("F146").Value = 'About right'
$ws.Range("G146").Value = 'A good broad introduction to the different systems at play and understanding the terminology'

$ws.Range("A147").Value = 45474.58373842593
$ws.Range("B147").Value = 45474.5841087963
$ws.Range("C147").Value = 'An introduction to AI'
$ws.Range("D147").Value = 45474
$ws.Range("E147").Value = 'Yes'
$ws.Range("F147").Value = 'About right'
$ws.Range("G147").Value = 'very interesting'

$ws.Range("A148").Value = 45474.58342592593
$ws.Range("B148").Value = 45474.584444444445
$ws.Range("C148").Value = 'An introduction to AI'
$ws.Range("D148").Value = 45474
$ws.Range("E148").Value = 'Yes'
$ws.Range("F148").Value = 'About right'
$ws.Range("G148").Value = 'As a clinician who is interested in data, AI and so on, I learned a lot and enjoyed the whole session.  '

$ws.Range("A149").Value = 45474.583657407406
$ws.Range("B149").Value = 45474.584641203706
$ws.Range("C149").Value = 'An introduction to AI'
$ws.Range("D149").Value = 45474
$ws.Range("E149").Value = 'Yes'
$ws.Range("F149").Value = 'About right'
$ws.Range("G149").Value = 'Very nice introduction and exciting start for staff to familiarise ourselves with AI concepts.'

$ws.Range("A150").Value = 45474.58332175926
$ws.Range("B150").Value = 45474.584652777776
$ws.Range("C150").Value = 'An Introduction to AI'
$ws.Range("D150").Value = 45474
$ws.Range("E150").Value = 'Yes'
$ws.Range("F150").Value = 'About right'

$ws.Range("A151").Value = 45474.58335648148
$ws.Range("B151").Value = 45474.584652777776
$ws.Range("C151").Value = 'An introduction to AI'
$ws.Range("D151").Value = 45474
$ws.Range("E151").Value = 'Yes'
$ws.Range("F151").Value = 'About right'
$ws.Range("G151").Value = 'Really helpful overview which prompted lots of discussion in the chat. '

$ws.Range("A152").Value = 45474.583449074074
$ws.Range("B152").Value = 45474.584710648145
$ws.Range("C152").Value = 'An introduction to AI'
$ws.Range("D152").Value = 45474
$ws.Range("E152").Value = 'Yes'
$ws.Range("F152").Value = 'About right'
$ws.Range("G152").Value = 'Good basic info but flitted around a lot on things I didn''t know - felt a bit unfocused.  Needed a bit more structure or someone to manage the chat that wasn''t the presenter.'

$ws.Range("A153").Value = 45474.58380787037
$ws.Range("B153").Value = 45474.5847337963
$ws.Range("C153").Value = 'An introduction to AI'
$ws.Range("D153").Value = 45474
$ws.Range("E153").Value = 'Yes'
$ws.Range("F153").Value = 'About right'
$ws.Range("G153").Value = 'Good introduction to AI, the development of it, current uses and shortcomings, and possibilities for the future'

$ws.Range("A154").Value = 45474.58349537037
$ws.Range("B154").Value = 45474.5850462963
$ws.Range("C154").Value = 'An introduction to AI'
$ws.Range("D154").Value = 45474
$ws.Range("E154").Value = 'Yes'
$ws.Range("F154").Value = 'About right'
$ws.Range("G154").Value = 'the essentials were there, whatever the care and social health focus; some of the information tidied up a few thoughts I had about AI'

$ws.Range("A155").Value = 45474.58366898148
$ws.Range("B155").Value = 45474.585648148146
$ws.Range("C155").Value = 'An introduction to AI'
$ws.Range("D155").Value = 45474
$ws.Range("E155").Value = 'Yes'
$ws.Range("F155").Value = 'About right'
$ws.Range("G155").Value = 'Thought it was a good intro but would be interesting to look at the main differences between "old AI" and new generative AI, The session you did looking at expert systems was good - are these classed as AI? When I was younger they were but now seems to refer mainly to generative AI so interesting to look at some of the history/hype vs the reality and the fact that most of the new AI is black box, so you can''t trace the algorithm back to why it made a particular descision in the way that you could in an expert system. This makes reviewing the Ai''s decisions a much harder task..'

$ws.Range("A156").Value = 45474.58346064815
$ws.Range("B156").Value = 45474.586064814815
$ws.Range("C156").Value = 'An introduction to AI'
$ws.Range("D156").Value = 45474
$ws.Range("E156").Value = 'Yes'
$ws.Range("F156").Value = 'About right'
$ws.Range("G156").Value = 'A nice introduction to an incredibly complex topic. It is challenging and covers a lot of ground in a short time but learning new terms and ideas, and being forced to think about things in a different way is never a bad thing.  Some great resources linked into the session that I will definitely need to go away and read (as well some brilliant stuff in the chat).'

$ws.Range("A157").Value = 45474.58393518518
$ws.Range("B157").Value = 45474.5878125
$ws.Range("C157").Value = 'An introduction to AI'
$ws.Range("D157").Value = 45474
$ws.Range("E157").Value = 'Yes'
$ws.Range("F157").Value = 'About right'
$ws.Range("G157").Value = 'Great overview, fascinating discussion - I really enjoy all of the references to explore at KIND sessions'

$ws.Range("A158").Value = 45474.61523148148
$ws.Range("B158").Value = 45474.617210648146
$ws.Range("C158").Value = 'An introduction to AI'
$ws.Range("D158").Value = 45474
$ws.Range("E158").Value = 'Yes'
$ws.Range("F158").Value = 'About right'
$ws.Range("G158").Value = 'The side discussions are always an exciting and valuable part of the informal training sessions - even the ones which are only tenuously linked. The content could have been offered as a slide deck ... but it is the side discussions which add substantial value (and interest) beyond anything which could simply be trawled from the internet'

$ws.Range("A159").Value = 45476.66207175926
$ws.Range("B159").Value = 45476.66236111111
$ws.Range("C159").Value = 'formulas'
$ws.Range("D159").Value = 45476
$ws.Range("E159").Value = 'Yes'
$ws.Range("F159").Value = 'About right'

$ws.Range("A160").Value = 45476.66211805555
$ws.Range("B160").Value = 45476.66273148148
$ws.Range("C160").Value = 'Excel formulas'
$ws.Range("D160").Value = 45476
$ws.Range("E160").Value = 'Yes'
$ws.Range("F160").Value = 'About right'
$ws.Range("G160").Value = 'These are just brilliant, a version of which I write at the end of every session. Go to them, just go. If you''re not sure whether to go. Just go.'

$ws.Range("A161").Value = 45476.6622337963
$ws.Range("B161").Value = 45476.66275462963
$ws.Range("C161").Value = 'Excel formulas'
$ws.Range("D161").Value = 45476
$ws.Range("E161").Value = 'Yes'
$ws.Range("F161").Value = 'About right'
$ws.Range("G161").Value = 'Good session for an absolute beginner'

$ws.Range("A162").Value = 45476.66216435185
$ws.Range("B162").Value = 45476.66234953704
$ws.Range("C162").Value = 'Excel formulas'
$ws.Range("D162").Value = 45476
$ws.Range("E162").Value = 'Yes'
$ws.Range("F162").Value = 'About right'
$ws.Range("G162").Value = 'Really good beginner session'

$ws.Range("A163").Value = 45476.66207175926
$ws.Range("B163").Value = 45476.6625
$ws.Range("C163").Value = 'Excel formulas'
$ws.Range("D163").Value = 45476
$ws.Range("E163").Value = 'Yes'
$ws.Range("F163").Value = 'About right'
$ws.Range("G163").Value = 'Great session, easy to follow. Thanks Brendan'

$ws.Range("A164").Value = 45476.66261574074
$ws.Range("B164").Value = 45476.66290509259
$ws.Range("C164").Value = 'Excel formulas'
$ws.Range("D164").Value = 45476
$ws.Range("E164").Value = 'Yes'
$ws.Range("F164").Value = 'About right'
$ws.Range("G164").Value = 'short sharp straight to the point'

$ws.Range("A165").Value = 45476.662141203706
$ws.Range("B165").Value = 45476.662881944445
$ws.Range("C165").Value = 'Excel formulas'
$ws.Range("D165").Value = 45476
$ws.Range("E165").Value = 'Yes'
$ws.Range("F165").Value = 'About right'
$ws.Range("G165").Value = 'It is really helpful to be shown a task, and give us the chance to try it for ourselves. '

$ws.Range("A166").Value = 45476.662094907406
$ws.Range("B166").Value = 45476.662997685184
$ws.Range("C166").Value = 'Excel formulas'
$ws.Range("D166").Value = 45476
$ws.Range("E166").Value = 'Yes'
$ws.Range("F166").Value = 'About right'
$ws.Range("G166").Value = 'Great overview of formulas and I learnt about the formula builder within Excel, that as a self-taught user I''d never used before which was great!'

$ws.Range("A167").Value = 45476.662083333336
$ws.Range("B167").Value = 45476.66304398148
$ws.Range("C167").Value = 'Excel Formulas'
$ws.Range("D167").Value = 45476
$ws.Range("E167").Value = 'Yes'
$ws.Range("F167").Value = 'About right'
$ws.Range("G167").Value = 'Enjoyed the session '

$ws.Range("A168").Value = 45476.662314814814
$ws.Range("B168").Value = 45476.66307870371
$ws.Range("C168").Value = 'Excel Formulas'
$ws.Range("D168").Value = 45476
$ws.Range("E168").Value = 'Yes'
$ws.Range("F168").Value = 'About right'
$ws.Range("G168").Value = 'It was informative and helpful.'

$ws.Range("A169").Value = 45476.66206018518
$ws.Range("B169").Value = 45476.66321759259
$ws.Range("C169").Value = 'excel formatting'
$ws.Range("D169").Value = 45476
$ws.Range("E169").Value = 'Yes'
$ws.Range("F169").Value = 'About right'
$ws.Range("G169").Value = 'Session was good and as described.  More focus on specific formulas and the reasons we may use them could potentially be beneficial.'

$ws.Range("C170").Value = 'Excel formulas'
$ws.Range("D170").Value = 45476
$ws.Range("E170").Value = 'Yes'
$ws.Range("F170").Value = 'About right'
$ws.Range("G170").Value = 'Very clear, helpful practical intro / review of formulas and simple functions in Excel'

$ws.Range("C171").Value = 'Excel formulas'
$ws.Range("D171").Value = 45476
$ws.Range("E171").Value = 'Yes'
$ws.Range("F171").Value = 'About right'
$ws.Range("G171").Value = 'Another great session - good balance of refreshing knowledge and learning new functions I wasn''t aware of!'

$ws.Range("A172").Value = 45476.66226851852
$ws.Range("B172").Value = 45476.66421296296
$ws.Range("C172").Value = 'Excel formulas'
$ws.Range("D172").Value = 45476
$ws.Range("E172").Value = 'Yes'
$ws.Range("F172").Value = 'About right'
$ws.Range("G172").Value = 'Was well done and think its important to provide training for excel as it probably the most under utilised software in our field.'

$ws.Range("C173").Value = 'Excel formulas'
$ws.Range("D173").Value = 45476
$ws.Range("E173").Value = 'Yes'
$ws.Range("F173").Value = 'About right'
$ws.Range("G173").Value = 'Basic beginner level - you could probably go on this same course a few times to try and get the basics, building on your knowledge each time.  You may pick up different things each session so I''ll try and attend again the next time this course is run.'

$ws.Range("A174").Value = 45476.66212962963
$ws.Range("B174").Value = 45476.66527777778
$ws.Range("C174").Value = 'Excel formulas'
$ws.Range("D174").Value = 45476
$ws.Range("E174").Value = 'Yes'
$ws.Range("F174").Value = 'About right'
$ws.Range("G174").Value = 'Really useful for going over basic concepts that I haven''t studied in a long time. So it was great for relearning the basic building blocks of how Excel works. Appropriate level for beginners session.'

$ws.Range("A175").Value = 45476.66271990741
$ws.Range("B175").Value = 45476.6662037037
$ws.Range("C175").Value = 'Kind Learning Network - Formulas in Excel'
$ws.Range("D175").Value = 45476
$ws.Range("E175").Value = 'Yes'
$ws.Range("F175").Value = 'About right'
$ws.Range("G175").Value = 'I really enjoyed the session, I had some (minimal) previous experience of Excel.  I found the examples were demonstrated well and easy to try during the session.  I couldn''t access the sample data as needed a sign in, I managed to copy and paste it over to the excel I was using but missed how to add the extra column.  It was great to be able to work out the average, will test that on some data I have collected alongside the other tips provided.  Many thanks for my place on the session.'

$ws.Range("A176").Value = 45476.66502314815
$ws.Range("B176").Value = 45476.66684027778
$ws.Range("C176").Value = 'Formulas in Excel '
$ws.Range("D176").Value = 45476
$ws.Range("E176").Value = 'Yes'
$ws.Range("F176").Value = 'About right'
$ws.Range("G176").Value = 'Easy to follow.  Clear instructions and good communication to provide session material prior to training. '

$ws.Range("A177").Value = 45476.66217592593
$ws.Range("B177").Value = 45477.63611111111
$ws.Range("C177").Value = 'Excel formulas'
$ws.Range("D177").Value = 45476
$ws.Range("E177").Value = 'Yes'
$ws.Range("F177").Value = 'About right'
$ws.Range("G177").Value = 'Some helpful tips and tricks, even for folks who were aware of formulae already.'

$ws.Range("A178").Value = 45478.560648148145
$ws.Range("B178").Value = 45478.56190972222
$ws.Range("C178").Value = 'Iteration in R'
$ws.Range("D178").Value = 45478
$ws.Range("E178").Value = 'Yes'
$ws.Range("F178").Value = 'About right'
$ws.Range("G178").Value = 'Excellent introduction to iteration - easy to follow along during the session'

$ws.Range("A179").Value = 45478.55978009259
$ws.Range("B179").Value = 45478.56466435185
$ws.Range("C179").Value = 'Iteration in R'
$ws.Range("D179").Value = 45478
$ws.Range("E179").Value = 'Yes'
$ws.Range("F179").Value = 'About right'
$ws.Range("G179").Value = 'Great session on For Loops and how to use them, also discussed why they potentially are not used as much as in i.e. Python'

$ws.Range("A180").Value = 45478.56922453704
$ws.Range("B180").Value = 45478.570185185185
$ws.Range("C180").Value = 'Shiny from scratch'
$ws.Range("D180").Value = 45476
$ws.Range("E180").Value = 'Yes'
$ws.Range("F180").Value = 'About right'
$ws.Range("G180").Value = 'Thought it was a great introduction to shiny, covered the core concepts well and gave me enough of an understanding that i''d feel confident going away and trying to build my own dashboard. '

$ws.Range("A181").Value = 45478.56972222222
$ws.Range("B181").Value = 45478.570231481484
$ws.Range("C181").Value = 'Shiny from scratch'
$ws.Range("D181").Value = 45476
$ws.Range("E181").Value = 'Yes'
$ws.Range("F181").Value = 'About right'
$ws.Range("G181").Value = 'It was great and well paced. Easy to follow along and code.'

$ws.Range("A182").Value = 45478.56914351852
$ws.Range("B182").Value = 45478.57109953704
$ws.Range("C182").Value = 'Shiny from scratch'
$ws.Range("D182").Value = 45476
$ws.Range("E182").Value = 'Yes'
$ws.Range("F182").Value = 'About right'
$ws.Range("G182").Value = 'Great introduction to the basics - didn''t know where to start with Shiny, and now I do!'

$ws.Range("A183").Value = 45478.57048611111
$ws.Range("B183").Value = 45478.5712037037
$ws.Range("C183").Value = 'Iteration in R'
$ws.Range("D183").Value = 45478
$ws.Range("E183").Value = 'Yes'
$ws.Range("F183").Value = 'About right'
$ws.Range("G183").Value = 'good introduction to loops - could maybe have had some more complex examples'

$ws.Range("A184").Value = 45478.574016203704
$ws.Range("B184").Value = 45478.57436342593
$ws.Range("C184").Value = 'Iteration in R'
$ws.Range("D184").Value = 45478
$ws.Range("E184").Value = 'Yes'
$ws.Range("F184").Value = 'About right'

$ws.Range("A185").Value = 45478.576689814814
$ws.Range("B185").Value = 45478.57722222222
$ws.Range("C185").Value = 'Shiny from scratch'
$ws.Range("D185").Value = 45476
$ws.Range("E185").Value = 'Yes'
$ws.Range("F185").Value = 'About right'
$ws.Range("G185").Value = 'Very nice introduction to the basics of Shiny!'

$ws.Range("A186").Value = 45478.580243055556
$ws.Range("B186").Value = 45478.58043981482
$ws.Range("C186").Value = 'Iteration in R'
$ws.Range("D186").Value = 45478
$ws.Range("E186").Value = 'No'
$ws.Range("F186").Value = 'Too easy / simple relative to description'

$ws.Range("A187").Value = 45478.58559027778
$ws.Range("B187").Value = 45478.58861111111
$ws.Range("C187").Value = 'Iteration in R'
$ws.Range("D187").Value = 45478
$ws.Range("E187").Value = 'Yes'
$ws.Range("F187").Value = 'About right'
$ws.Range("G187").Value = 'A great into to iteration!'

$ws.Range("A188").Value = 45478.5883912037
$ws.Range("B188").Value = 45478.58914351852
$ws.Range("C188").Value = 'Iteration in R'
$ws.Range("D188").Value = 45478
$ws.Range("E188").Value = 'Yes'
$ws.Range("F188").Value = 'About right'

$ws.Range("A189").Value = 45478.5890162037
$ws.Range("B189").Value = 45478.58954861111
$ws.Range("C189").Value = 'Shiny from scratch'
$ws.Range("D189").Value = 45476
$ws.Range("E189").Value = 'Yes'
$ws.Range("F189").Value = 'About right'
$ws.Range("G189").Value = 'Great intro session to Shiny, which made me want to learn more!'

$ws.Range("A190").Value = 45478.58974537037
$ws.Range("B190").Value = 45478.58996527778
$ws.Range("C190").Value = 'Shiny from scratch'
$ws.Range("D190").Value = 45476
$ws.Range("E190").Value = 'Yes'
$ws.Range("F190").Value = 'About right'
$ws.Range("G190").Value = 'Great starter course for Shiny'

# Resize the table/list object to cover the new data range
$lo.Resize($ws.Range("A1:G190"))

# Update the hidden ExternalData_1 defined name to match the new table range
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Table1!ExternalData_1") {
        $n.RefersTo = "=Table1!`$A`$1:`$G`$190"
    }
}

# Update the view: scroll position and active cell selection
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 166
$ws.Range("C186").Select()

